$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking data columns (G:K) to Text format so values
# like "0", "0.00" are stored as text, matching the source data.
$ws.Range("G2:K18").NumberFormat = "@"

# Header row (row 1)
$ws.Range("A1").Value = 'venue'
$ws.Range("B1").Value = 'date'
$ws.Range("C1").Value = 'result'
$ws.Range("D1").Value = 'ownTeam'
$ws.Range("E1").Value = 'oppTeam'
$ws.Range("F1").Value = 'batsman'
$ws.Range("G1").Value = 'totalRuns'
$ws.Range("H1").Value = 'totalBalls'
$ws.Range("I1").Value = 'total4s'
$ws.Range("J1").Value = 'total6s'
$ws.Range("K1").Value = 'sr'

# Data rows (2-18)
# Row 2
$ws.Range("A2").Value = ' Dubai (DSC)'
$ws.Range("B2").Value = ' October 27 2020'
$ws.Range("C2").Value = 'Sunrisers won by 88 runs'
$ws.Range("D2").Value = 'Delhi Capitals'
$ws.Range("E2").Value = 'Sunrisers Hyderabad'
$ws.Range("F2").Value = 'Shikhar Dhawan '
$ws.Range("G2").Value = '0'
$ws.Range("H2").Value = '1'
$ws.Range("I2").Value = '0'
$ws.Range("J2").Value = '0'
$ws.Range("K2").Value = '0.00'

# Row 3
$ws.Range("A3").Value = ' Dubai (DSC)'
$ws.Range("B3").Value = ' November 05 2020'
$ws.Range("C3").Value = 'Mumbai won by 57 runs'
$ws.Range("D3").Value = 'Delhi Capitals'
$ws.Range("E3").Value = 'Mumbai Indians'
$ws.Range("F3").Value = 'Shikhar Dhawan '
$ws.Range("G3").Value = '0'
$ws.Range("H3").Value = '2'
$ws.Range("I3").Value = '0'
$ws.Range("J3").Value = '0'
$ws.Range("K3").Value = '0.00'

# Row 4
$ws.Range("A4").Value = ' Dubai (DSC)'
$ws.Range("B4").Value = ' November 10 2020'
$ws.Range("C4").Value = 'Mumbai won by 5 wickets (with 8 balls remaining)'
$ws.Range("D4").Value = 'Delhi Capitals'
$ws.Range("E4").Value = 'Mumbai Indians'
$ws.Range("F4").Value = 'Shikhar Dhawan '
$ws.Range("G4").Value = '15'
$ws.Range("H4").Value = '13'
$ws.Range("I4").Value = '3'
$ws.Range("J4").Value = '0'
$ws.Range("K4").Value = '115.38'

# Row 5
$ws.Range("A5").Value = ' Dubai (DSC)'
$ws.Range("B5").Value = ' October 14 2020'
$ws.Range("C5").Value = 'Capitals won by 13 runs'
$ws.Range("D5").Value = 'Delhi Capitals'
$ws.Range("E5").Value = 'Rajasthan Royals'
$ws.Range("F5").Value = 'Shikhar Dhawan '
$ws.Range("G5").Value = '57'
$ws.Range("H5").Value = '33'
$ws.Range("I5").Value = '6'
$ws.Range("J5").Value = '2'
$ws.Range("K5").Value = '172.72'

# Row 6
$ws.Range("A6").Value = ' Dubai (DSC)'
$ws.Range("B6").Value = ' October 31 2020'
$ws.Range("C6").Value = 'Mumbai won by 9 wickets (with 34 balls remaining)'
$ws.Range("D6").Value = 'Delhi Capitals'
$ws.Range("E6").Value = 'Mumbai Indians'
$ws.Range("F6").Value = 'Shikhar Dhawan '
$ws.Range("G6").Value = '0'
$ws.Range("H6").Value = '2'
$ws.Range("I6").Value = '0'
$ws.Range("J6").Value = '0'
$ws.Range("K6").Value = '0.00'

# Row 7
$ws.Range("A7").Value = ' Dubai (DSC)'
$ws.Range("B7").Value = ' September 25 2020'
$ws.Range("C7").Value = 'Capitals won by 44 runs'
$ws.Range("D7").Value = 'Delhi Capitals'
$ws.Range("E7").Value = 'Chennai Super Kings'
$ws.Range("F7").Value = 'Shikhar Dhawan '
$ws.Range("G7").Value = '35'
$ws.Range("H7").Value = '27'
$ws.Range("I7").Value = '3'
$ws.Range("J7").Value = '1'
$ws.Range("K7").Value = '129.62'

# Row 8
$ws.Range("A8").Value = ' Sharjah'
$ws.Range("B8").Value = ' October 17 2020'
$ws.Range("C8").Value = 'Capitals won by 5 wickets (with 1 ball remaining)'
$ws.Range("D8").Value = 'Delhi Capitals'
$ws.Range("E8").Value = 'Chennai Super Kings'
$ws.Range("F8").Value = 'Shikhar Dhawan '
$ws.Range("G8").Value = '101'
$ws.Range("H8").Value = '58'
$ws.Range("I8").Value = '14'
$ws.Range("J8").Value = '1'
$ws.Range("K8").Value = '174.13'

# Row 9
$ws.Range("A9").Value = ' Sharjah'
$ws.Range("B9").Value = ' October 09 2020'
$ws.Range("C9").Value = 'Capitals won by 46 runs'
$ws.Range("D9").Value = 'Delhi Capitals'
$ws.Range("E9").Value = 'Rajasthan Royals'
$ws.Range("F9").Value = 'Shikhar Dhawan '
$ws.Range("G9").Value = '5'
$ws.Range("H9").Value = '4'
$ws.Range("I9").Value = '1'
$ws.Range("J9").Value = '0'
$ws.Range("K9").Value = '125.00'

# Row 10
$ws.Range("A10").Value = ' Dubai (DSC)'
$ws.Range("B10").Value = ' October 05 2020'
$ws.Range("C10").Value = 'Capitals won by 59 runs'
$ws.Range("D10").Value = 'Delhi Capitals'
$ws.Range("E10").Value = 'Royal Challengers Bangalore'
$ws.Range("F10").Value = 'Shikhar Dhawan '
$ws.Range("G10").Value = '32'
$ws.Range("H10").Value = '28'
$ws.Range("I10").Value = '3'
$ws.Range("J10").Value = '0'
$ws.Range("K10").Value = '114.28'

# Row 11
$ws.Range("A11").Value = ' Dubai (DSC)'
$ws.Range("B11").Value = ' October 20 2020'
$ws.Range("C11").Value = 'Kings XI won by 5 wickets (with 6 balls remaining)'
$ws.Range("D11").Value = 'Delhi Capitals'
$ws.Range("E11").Value = 'Kings XI Punjab'
$ws.Range("F11").Value = 'Shikhar Dhawan '
$ws.Range("G11").Value = '106'
$ws.Range("H11").Value = '61'
$ws.Range("I11").Value = '12'
$ws.Range("J11").Value = '3'
$ws.Range("K11").Value = '173.77'

# Row 12
$ws.Range("A12").Value = ' Dubai (DSC)'
$ws.Range("B12").Value = ' September 20 2020'
$ws.Range("C12").Value = 'Match tied (Capitals won the one-over eliminator)'
$ws.Range("D12").Value = 'Delhi Capitals'
$ws.Range("E12").Value = 'Kings XI Punjab'
$ws.Range("F12").Value = 'Shikhar Dhawan '
$ws.Range("G12").Value = '0'
$ws.Range("H12").Value = '2'
$ws.Range("I12").Value = '0'
$ws.Range("J12").Value = '0'
$ws.Range("K12").Value = '0.00'

# Row 13
$ws.Range("A13").Value = ' Sharjah'
$ws.Range("B13").Value = ' October 03 2020'
$ws.Range("C13").Value = 'Capitals won by 18 runs'
$ws.Range("D13").Value = 'Delhi Capitals'
$ws.Range("E13").Value = 'Kolkata Knight Riders'
$ws.Range("F13").Value = 'Shikhar Dhawan '
$ws.Range("G13").Value = '26'
$ws.Range("H13").Value = '16'
$ws.Range("I13").Value = '2'
$ws.Range("J13").Value = '2'
$ws.Range("K13").Value = '162.50'

# Row 14
$ws.Range("A14").Value = ' Abu Dhabi'
$ws.Range("B14").Value = ' November 08 2020'
$ws.Range("C14").Value = 'Capitals won by 17 runs'
$ws.Range("D14").Value = 'Delhi Capitals'
$ws.Range("E14").Value = 'Sunrisers Hyderabad'
$ws.Range("F14").Value = 'Shikhar Dhawan '
$ws.Range("G14").Value = '78'
$ws.Range("H14").Value = '50'
$ws.Range("I14").Value = '6'
$ws.Range("J14").Value = '2'
$ws.Range("K14").Value = '156.00'

# Row 15
$ws.Range("A15").Value = ' Abu Dhabi'
$ws.Range("B15").Value = ' November 02 2020'
$ws.Range("C15").Value = 'Capitals won by 6 wickets (with 6 balls remaining)'
$ws.Range("D15").Value = 'Delhi Capitals'
$ws.Range("E15").Value = 'Royal Challengers Bangalore'
$ws.Range("F15").Value = 'Shikhar Dhawan '
$ws.Range("G15").Value = '54'
$ws.Range("H15").Value = '41'
$ws.Range("I15").Value = '6'
$ws.Range("J15").Value = '0'
$ws.Range("K15").Value = '131.70'

# Row 16
$ws.Range("A16").Value = ' Abu Dhabi'
$ws.Range("B16").Value = ' October 24 2020'
$ws.Range("C16").Value = 'KKR won by 59 runs'
$ws.Range("D16").Value = 'Delhi Capitals'
$ws.Range("E16").Value = 'Kolkata Knight Riders'
$ws.Range("F16").Value = 'Shikhar Dhawan '
$ws.Range("G16").Value = '6'
$ws.Range("H16").Value = '6'
$ws.Range("I16").Value = '1'
$ws.Range("J16").Value = '0'
$ws.Range("K16").Value = '100.00'

# Row 17
$ws.Range("A17").Value = ' Abu Dhabi'
$ws.Range("B17").Value = ' October 11 2020'
$ws.Range("C17").Value = 'Mumbai won by 5 wickets (with 2 balls remaining)'
$ws.Range("D17").Value = 'Delhi Capitals'
$ws.Range("E17").Value = 'Mumbai Indians'
$ws.Range("F17").Value = 'Shikhar Dhawan '
$ws.Range("G17").Value = '69'
$ws.Range("H17").Value = '52'
$ws.Range("I17").Value = '6'
$ws.Range("J17").Value = '1'
$ws.Range("K17").Value = '132.69'

# Row 18
$ws.Range("A18").Value = ' Abu Dhabi'
$ws.Range("B18").Value = ' September 29 2020'
$ws.Range("C18").Value = 'Sunrisers won by 15 runs'
$ws.Range("D18").Value = 'Delhi Capitals'
$ws.Range("E18").Value = 'Sunrisers Hyderabad'
$ws.Range("F18").Value = 'Shikhar Dhawan '
$ws.Range("G18").Value = '34'
$ws.Range("H18").Value = '31'
$ws.Range("I18").Value = '4'
$ws.Range("J18").Value = '0'
$ws.Range("K18").Value = '109.67'
